$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the KPI result values (row 2 - row 5, columns B/C/D)
$ws.Range("B2").Value = 83754.347093073899
$ws.Range("C2").Value = 60269.168400622002
$ws.Range("D2").Value = 373.30580949705802

$ws.Range("B3").Value = 82214.984532106304
$ws.Range("C3").Value = 58775.838255845403
$ws.Range("D3").Value = 419.33822568828901

$ws.Range("B4").Value = 80743.654556177498
$ws.Range("C4").Value = 57316.646779459297
$ws.Range("D4").Value = 431.47672523089102

$ws.Range("B5").Value = 78348.550411831704
$ws.Range("C5").Value = 54993.302038123802
$ws.Range("D5").Value = 503.23612824068499

# Update the active selection (was D6, now B9)
$excel.Goto($ws.Range("B9"))
